$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginDetails")

# --- Add new worksheets, each placed after the previous one, so final
#     order is: LoginDetails, AdminUserDetails, CategoryDetail,
#     ContactDetails, FooterDetails, NewDetails, SubCategoryDetails ---
$wsAdmin = $wb.Worksheets.Add($null, $ws1)
$wsAdmin.Name = "AdminUserDetails"

$wsCategory = $wb.Worksheets.Add($null, $wsAdmin)
$wsCategory.Name = "CategoryDetail"

$wsContact = $wb.Worksheets.Add($null, $wsCategory)
$wsContact.Name = "ContactDetails"

$wsFooter = $wb.Worksheets.Add($null, $wsContact)
$wsFooter.Name = "FooterDetails"

$wsNews = $wb.Worksheets.Add($null, $wsFooter)
$wsNews.Name = "NewDetails"

$wsSubCategory = $wb.Worksheets.Add($null, $wsNews)
$wsSubCategory.Name = "SubCategoryDetails"

# --- AdminUserDetails ---
$wsAdmin.Range("A1").Value = "UserName"
$wsAdmin.Range("B1").Value = "Password"
$wsAdmin.Range("C1").Value = "AdminUserName"
$wsAdmin.Range("D1").Value = "AdminPassword"
$wsAdmin.Range("A2").Value = "admin"
$wsAdmin.Range("B2").Value = "admin"
$wsAdmin.Range("C2").Value = "Akhil"
$wsAdmin.Range("D2").Value = "User123"

# --- CategoryDetail ---
$wsCategory.Range("A1").Value = "UserName"
$wsCategory.Range("B1").Value = "Password"
$wsCategory.Range("C1").Value = "CategoryName"
$wsCategory.Range("A2").Value = "admin"
$wsCategory.Range("B2").Value = "admin"
$wsCategory.Range("C2").Value = "Test Category 1122"

# --- ContactDetails ---
$wsContact.Range("A1").Value = "UserName"
$wsContact.Range("B1").Value = "Password"
$wsContact.Range("C1").Value = "PhoneNo"
$wsContact.Range("D1").Value = "Email"
$wsContact.Range("E1").Value = "Address"
$wsContact.Range("F1").Value = "Delivery Time"
$wsContact.Range("G1").Value = "Delivery Charge"
$wsContact.Range("A2").Value = "admin"
$wsContact.Range("B2").Value = "admin"
$wsContact.Range("C2").Value = 9995870445
$wsContact.Range("D2").Value = "test1122@mail.com"
$wsContact.Hyperlinks.Add($wsContact.Range("D2"), "mailto:test1122@mail.com") | Out-Null
$wsContact.Range("E2").Value = "Phase 1, Technopark"
$wsContact.Range("F2").Value = 10
$wsContact.Range("G2").Value = 40

# --- FooterDetails ---
$wsFooter.Range("A1").Value = "UserName"
$wsFooter.Range("B1").Value = "Password"
$wsFooter.Range("C1").Value = "Address"
$wsFooter.Range("D1").Value = "Email"
$wsFooter.Range("E1").Value = "Phone"
$wsFooter.Range("A2").Value = "admin"
$wsFooter.Range("B2").Value = "admin"
$wsFooter.Range("C2").Value = "Asiatic business center,Technopark Phase"
$wsFooter.Range("D2").Value = "testing@gmail.com"
$wsFooter.Hyperlinks.Add($wsFooter.Range("D2"), "mailto:testing@gmail.com") | Out-Null
# Phone number stored as text (quote-prefixed) so it isn't treated as a number
$wsFooter.Range("E2").Value = "'8947584758"

# --- NewDetails ---
$wsNews.Range("A1").Value = "UserName"
$wsNews.Range("B1").Value = "Password"
$wsNews.Range("C1").Value = "News"
$wsNews.Range("A2").Value = "admin"
$wsNews.Range("B2").Value = "admin"
$wsNews.Range("C2").Value = "Lorem Ipsum is simply dummy text of the printing and typesetting industry."

# --- SubCategoryDetails ---
$wsSubCategory.Range("A1").Value = "UserName"
$wsSubCategory.Range("B1").Value = "Password"
$wsSubCategory.Range("C1").Value = "SubCategoryName"
$wsSubCategory.Range("A2").Value = "admin"
$wsSubCategory.Range("B2").Value = "admin"
$wsSubCategory.Range("C2").Value = "Test Subcategory 1126"

# --- Selections to match target output ---
$ws1.Range("A1:B1").Select()
$wsAdmin.Range("A1:B2").Select()
$wsCategory.Range("A1:B2").Select()
$wsContact.Range("A1:B2").Select()
$wsFooter.Range("A1:B2").Select()
$wsNews.Range("A1:B2").Select()
$wsSubCategory.Range("J19").Select()

$wsSubCategory.Activate()
$excel.Windows.Item(1).ScrollWorkbookTabs(1)
